$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restyle rows 35 and 36 to match the "closed" bordered look used by the
# rows above them (copy cell formatting from row 34, which already carries
# the desired styles: date-style + border for column A, bottom/side border
# style for column B) ---
$ws.Range("A34").Copy()
$ws.Range("A35:A36").PasteSpecial(-4122)

$ws.Range("B34").Copy()
$ws.Range("B35:B36").PasteSpecial(-4122)

# Match the (taller) row heights that come with the new bordered look
$ws.Rows.Item(35).RowHeight = 35
$ws.Rows.Item(36).RowHeight = 18

# --- Add the new task row 37 (date + task name), reusing the same format ---
$ws.Range("A34").Copy()
$ws.Range("A37").PasteSpecial(-4122)
$ws.Range("A37").Value = 45184

$ws.Range("B34").Copy()
$ws.Range("B37").PasteSpecial(-4122)
$ws.Range("B37").Value = "fuzzy full-text search"

$ws.Rows.Item(37).RowHeight = 18

$excel.CutCopyMode = $false

# --- Update the view to reflect where the user ended up working ---
$ws.Range("B42").Select()
